$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells, matching style of existing header row (F1)
$ws.Range("F1").Copy()
$ws.Range("G1:H1").PasteSpecial(-4122)
$ws.Range("G1").Value = "Elapsed Time"
$ws.Range("H1").Value = "CPU"

# New data cells in row 2
$ws.Range("G2").Value = 0.1311458841167526
$ws.Range("H2").Value = 0.991
